# Apply the "update excel and upload source files" edit to the workbook.
# Sheets: "data" (risk survey results), "dates" (wave start/end/n), "source".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": fill in newly-collected wave columns (O, Q, R, S, T, AA,
# AB, AC) for several existing rows, and append a brand-new risk row 31.
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("data")

# Row 4 - Spannungen durch Zuzug von Auslaendern
$data.Range("O4").Value = 45
$data.Range("Q4").Value = 41
$data.Range("R4").Value = 37
$data.Range("S4").Value = 37
$data.Range("T4").Value = 42

# Row 7 - sinkender Lebensstandard im Alter
$data.Range("O7").Value = 48
$data.Range("Q7").Value = 38
$data.Range("R7").Value = 41
$data.Range("S7").Value = 37
$data.Range("T7").Value = 47
$data.Range("AA7").Value = 37
$data.Range("AB7").Value = 39
$data.Range("AC7").Value = 34

# Row 8 - Vereinsamung im Alter
$data.Range("O8").Value = 40
$data.Range("Q8").Value = 30
$data.Range("R8").Value = 31
$data.Range("S8").Value = 30
$data.Range("T8").Value = 36

# Row 10 - hoehere Arbeitslosigkeit in Deutschland
$data.Range("O10").Value = 48
$data.Range("Q10").Value = 43
$data.Range("R10").Value = 42
$data.Range("S10").Value = 34
$data.Range("T10").Value = 46
$data.Range("AA10").Value = 28
$data.Range("AB10").Value = 27
$data.Range("AC10").Value = 24

# Row 16 - Naturkatastrophen/Wetterextreme
$data.Range("O16").Value = 49
$data.Range("T16").Value = 64

# Row 20 - Ueberforderung des Staates durch Gefluechtete
$data.Range("AA20").Value = 42

# Row 29 - Zerbrechen der Partnerschaft
$data.Range("O29").Value = 24
$data.Range("P29").Value = 18
$data.Range("Q29").Value = 19
$data.Range("R29").Value = 21
$data.Range("S29").Value = 16
$data.Range("T29").Value = 23
$data.Range("AA29").Value = 17
$data.Range("AB29").Value = 18
$data.Range("AC29").Value = 18

# New row 31 - a brand-new risk factor
$data.Range("A31").Value = "fehlender Bürgernähe der Politiker"
$data.Range("O31").Value = 64
$data.Range("P31").Value = 63

# ---------------------------------------------------------------------
# Sheet "dates": backfill wave start/end dates for 2015 & 2016, and
# refresh the respondent counts (column E) for 2015, 2016 and 2020.
# ---------------------------------------------------------------------
$dates = $wb.Worksheets.Item("dates")

# Reuse the date number-format already applied to C28/D28 (and friends)
# so the new cells share the same style index instead of minting a new one.
$dates.Range("C28").Copy()
$dates.Range("C25").PasteSpecial(-4122)
$dates.Range("D25").PasteSpecial(-4122)
$dates.Range("C26").PasteSpecial(-4122)
$dates.Range("D26").PasteSpecial(-4122)

$dates.Range("C25").Value = 42160
$dates.Range("D25").Value = 42172
$dates.Range("E25").Value = 2373

$dates.Range("C26").Value = 42466
$dates.Range("D26").Value = 42503
$dates.Range("E26").Value = 2421

$dates.Range("E30").Value = 2396

# ---------------------------------------------------------------------
# Restore the current selections to match the saved-file state: the
# "dates" selection moves first (without activating it) then "data"
# becomes/stays the active sheet with its new selection.
# ---------------------------------------------------------------------
[void]$dates.Cells.Item(25, 5).Select()

[void]$data.Activate()
[void]$data.Cells.Item(31, 16).Select()
